# Applies the "苏州-漫展信息" update (commit: "Update gh-pages to output
# generated at 456a3b4") to both the "展览" and "全部类型" worksheets,
# which carry identical data tables.

function Set-TextCell($ws, $row, $col, $val) {
    # Force text storage so date-shaped strings like "2024-06-29" are not
    # auto-coerced into Excel date serials, then drop back to the default
    # "Normal" style so no stray number-format style sticks to the cell.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Set-NumCell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

function Update-SuzhouSheet($ws) {
    # --- simple "want-to-go" count bumps -----------------------------
    Set-NumCell $ws 5  6 12999
    Set-NumCell $ws 8  6 513
    Set-NumCell $ws 9  6 478
    Set-NumCell $ws 10 6 1171
    Set-NumCell $ws 11 6 973
    Set-NumCell $ws 12 6 13736
    Set-NumCell $ws 13 6 14275

    # --- rows 18-21: cv meet-and-greet rows move up one slot, and the
    # vacated slot (21) becomes the (renamed) main carnival listing ----
    Set-TextCell $ws 18 3 "苏州·动漫游戏嘉年华cv刘圣博见面会（取消）"
    Set-NumCell  $ws 18 6 1
    Set-TextCell $ws 18 7 "不可售"
    Set-TextCell $ws 18 8 "https://show.bilibili.com/platform/detail.html?id=83038"
    Set-TextCell $ws 18 9 "//i0.hdslb.com/bfs/openplatform/202403/D94B39u21710901393375.jpeg"

    Set-TextCell $ws 19 3 "苏州·动漫游戏嘉年华cv张文钰见面会（取消）"
    Set-NumCell  $ws 19 6 2
    Set-TextCell $ws 19 7 "不可售"
    Set-TextCell $ws 19 8 "https://show.bilibili.com/platform/detail.html?id=83037"
    Set-TextCell $ws 19 9 "//i1.hdslb.com/bfs/openplatform/202403/eP5thEuS1710901472569.jpeg"

    Set-TextCell $ws 20 3 "苏州·动漫游戏嘉年华cv沐霏见面会（取消）"
    Set-NumCell  $ws 20 6 4
    Set-TextCell $ws 20 7 "不可售"
    Set-TextCell $ws 20 8 "https://show.bilibili.com/platform/detail.html?id=82891"
    Set-TextCell $ws 20 9 "//i1.hdslb.com/bfs/openplatform/202403/8VORpvQc1710900704258.jpeg"

    Set-TextCell $ws 21 3 "苏州·苏州湾动漫游戏嘉年华"
    Set-NumCell  $ws 21 6 31
    Set-NumCell  $ws 21 7 58
    Set-TextCell $ws 21 8 "https://show.bilibili.com/platform/detail.html?id=82824"
    Set-TextCell $ws 21 9 "//i1.hdslb.com/bfs/openplatform/202403/HzWBEJeE1710324788092.jpeg"

    # --- more simple "want-to-go" count bumps -------------------------
    Set-NumCell $ws 25 6 5343
    Set-NumCell $ws 26 6 934
    Set-NumCell $ws 27 6 14
    Set-NumCell $ws 28 6 292

    # --- two brand-new listings appended at the bottom ----------------
    $ws.Cells.Item(28, 1).Copy()
    $ws.Cells.Item(29, 1).PasteSpecial(-4122)
    $ws.Cells.Item(30, 1).PasteSpecial(-4122)

    Set-NumCell  $ws 29 1 28
    Set-TextCell $ws 29 2 "2024-06-29"
    Set-TextCell $ws 29 3 "苏州·归离之缘原神only展"
    Set-TextCell $ws 29 4 "清禾路888号2号楼3楼 格莱美婚礼宴会中心"
    Set-TextCell $ws 29 5 "2024.06.29 09:30-06.29 18:30"
    Set-NumCell  $ws 29 6 6
    Set-NumCell  $ws 29 7 89
    Set-TextCell $ws 29 8 "https://show.bilibili.com/platform/detail.html?id=83271"
    Set-TextCell $ws 29 9 "//i1.hdslb.com/bfs/openplatform/202403/hNkSoRCt1710999968954.png"

    Set-NumCell  $ws 30 1 29
    Set-TextCell $ws 30 2 "2024-07-20"
    Set-TextCell $ws 30 3 "苏州·萤火国潮文化节动漫品牌博览会"
    Set-TextCell $ws 30 4 "金山南路288号木渎影视城F2 苏州广电国际会展中心"
    Set-TextCell $ws 30 5 "2024.07.20 10:00-07.21 17:00"
    Set-NumCell  $ws 30 6 0
    Set-NumCell  $ws 30 7 60
    Set-TextCell $ws 30 8 "https://show.bilibili.com/platform/detail.html?id=83301"
    Set-TextCell $ws 30 9 "//i0.hdslb.com/bfs/openplatform/202403/rV07luU61711274774556.jpeg"
}

$wb = $excel.ActiveWorkbook

Update-SuzhouSheet $wb.Worksheets.Item(1)   # 展览
Update-SuzhouSheet $wb.Worksheets.Item(4)   # 全部类型

Write-Output "edit complete"
